$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1, Q1 (style matches existing bold/bordered header style of row 1)
$headerRef = $ws.Range("O1")
foreach ($addr in @("P1", "Q1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $headerRef.Font.Bold
    $cell.Borders.LineStyle = $headerRef.Borders.LineStyle
    $cell.HorizontalAlignment = $headerRef.HorizontalAlignment
    $cell.VerticalAlignment = $headerRef.VerticalAlignment
}
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For data rows 2-25: swap values in columns I/K and M/O, then add new P/Q columns = 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column = 2
}
